# Update "想去人数" (want-to-go count) figures across sheets to match the
# latest scrape output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 30
$ws1.Range("F3").Value = 8950
$ws1.Range("F4").Value = 2701
$ws1.Range("F6").Value = 318
$ws1.Range("F7").Value = 832
$ws1.Range("F8").Value = 695
$ws1.Range("F9").Value = 129
$ws1.Range("F11").Value = 386
$ws1.Range("F13").Value = 3819
$ws1.Range("F14").Value = 288
$ws1.Range("F16").Value = 802
$ws1.Range("F19").Value = 501
$ws1.Range("F22").Value = 1385
$ws1.Range("F24").Value = 479
$ws1.Range("F27").Value = 173
$ws1.Range("F33").Value = 713
$ws1.Range("F36").Value = 98
$ws1.Range("F41").Value = 195
$ws1.Range("F42").Value = 328
$ws1.Range("F43").Value = 23
$ws1.Range("F44").Value = 15

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 2

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 30
$ws4.Range("F4").Value = 318
$ws4.Range("F5").Value = 832
$ws4.Range("F6").Value = 695
$ws4.Range("F7").Value = 129
$ws4.Range("F9").Value = 386
$ws4.Range("F12").Value = 3819
$ws4.Range("F13").Value = 288
$ws4.Range("F16").Value = 2
$ws4.Range("F17").Value = 802
$ws4.Range("F22").Value = 501
$ws4.Range("F26").Value = 1385
$ws4.Range("F28").Value = 479
$ws4.Range("F31").Value = 173
$ws4.Range("F37").Value = 713
$ws4.Range("F40").Value = 98
$ws4.Range("F44").Value = 195
$ws4.Range("F45").Value = 328
$ws4.Range("F46").Value = 23
$ws4.Range("F47").Value = 15
